$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update odds values in row 2 (columns M, N, O, P, S, T)
$ws.Range("M2").Value = 1.17
$ws.Range("N2").Value = 5
$ws.Range("O2").Value = 1.73
$ws.Range("P2").Value = 2
$ws.Range("S2").Value = 1.75
$ws.Range("T2").Value = 2.05

# Delete row 5 entirely (shifting cells up), reducing the used range to A1:BD4
$ws.Rows.Item(5).Delete()
